$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Prepend "color/" to the image file names in column A (rows 2-151),
#    e.g. "color_0.png" -> "color/color_0.png"
for ($r = 2; $r -le 151; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $old = $cell.Value2
    if ($old -like "color_*.png") {
        $cell.Value = "color/" + $old
    }
}

# 2) Update question text: "primary color" -> "main color" (trees question)
$ws.Range("B21").Value = "Describe the main color of the trees in the picture."

# 3) Update reference answer about butterfly colors to include "black"
$ws.Range("C23").Value = "Most of the butterflies in the picture are colored blue, light blue and black. There is a black and white butterfly in the picture. It is positioned near the center-left side among the flowers."
